$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Detaily")

# Update the unit price (B10) which drives the recalculation of F/G columns
$ws.Range("B10").Value = 24.51

# Activate the sheet and set the selection to A11 as in the saved file
$ws.Activate()
$ws.Range("A11").Select()
